# Daily attendance processing - 2026-01-04 11:01:01
# Reorders the "Last Modified By" list in column G: moves a leading
# "System" entry (the literal token "System", case-sensitive) from the
# front of the comma-separated list to the back of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value2

    if ($val -ne $null -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)   # strip off the leading "System, "
        $newVal = $rest + ", System"
        $cell.Value2 = $newVal
    }
}
